$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B4").Value = "inf"
$ws.Range("B6").Value = 57898.15688466049
$ws.Range("B7").Value = 10456577.88551109
$ws.Range("B8").Value = 24732210.50973683
$ws.Range("B10").Value = 2979694.134305789

# --- Fed-in Capacity sheet ---
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("R20").Value = 0
$ws.Range("J21").Value = 93.17061249236157
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 150.3839754851235
$ws.Range("I24").Value = 10.12574714858493
$ws.Range("J24").Value = 93.17061249236157
$ws.Range("K24").Value = 80.29914934735042
$ws.Range("L24").Value = 61.18167021676314
$ws.Range("O24").Value = 57.81213424001893
$ws.Range("Q24").Value = 94.49434172313325
$ws.Range("M25").Value = 92.09541281912071
$ws.Range("N25").Value = 81.96869489115805
$ws.Range("O25").Value = 96.22962838366004
$ws.Range("P25").Value = 101.5955875616828
$ws.Range("Q25").Value = 0
$ws.Range("Q29").Value = 150.3839754851235
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 65.92768427608706
$ws.Range("M32").Value = 0
$ws.Range("P32").Value = 135.4597561231036
$ws.Range("Q32").Value = 150.3839754851235
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 81.96869489115805
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 0
$ws.Range("R35").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("N39").Value = 0
$ws.Range("P39").Value = 0
$ws.Range("Q39").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = 0
$ws.Range("N40").Value = 0
$ws.Range("O40").Value = 0
$ws.Range("P40").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("P41").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = 0
$ws.Range("N43").Value = 0
$ws.Range("P43").Value = 101.5955875616828
$ws.Range("R44").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K46").Value = 94.30397654773019

# --- Unmet Demand sheet ---
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("L12").Value = 61.18167021676314
$ws.Range("M12").Value = 51.84373129681028
$ws.Range("O12").Value = 57.81213424001893
$ws.Range("M15").Value = 51.84373129681028
$ws.Range("L18").Value = 61.18167021676314
$ws.Range("N18").Value = 38.66169381481656
$ws.Range("L19").Value = 90.4687457914608
$ws.Range("P19").Value = 101.5955875616828
$ws.Range("R20").Value = 173.7492132756177
$ws.Range("J21").Value = 0
$ws.Range("P22").Value = 101.5955875616828
$ws.Range("Q22").Value = 126.4887893424616
$ws.Range("J23").Value = 124.5190384721106
$ws.Range("K23").Value = 135.370731907559
$ws.Range("P23").Value = 135.4597561231036
$ws.Range("Q23").Value = 0
$ws.Range("I24").Value = 77.12765456497084
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("Q24").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 126.4887893424616
$ws.Range("Q29").Value = 0
$ws.Range("L30").Value = 61.18167021676314
$ws.Range("M30").Value = 51.84373129681028
$ws.Range("N30").Value = 38.66169381481656
$ws.Range("O30").Value = 57.81213424001893
$ws.Range("P30").Value = 0
$ws.Range("M32").Value = 113.4004983079896
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("J33").Value = 93.17061249236157
$ws.Range("K33").Value = 80.29914934735042
$ws.Range("O33").Value = 57.81213424001893
$ws.Range("P33").Value = 65.92768427608706
$ws.Range("Q33").Value = 94.49434172313325
$ws.Range("M34").Value = 92.09541281912071
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 96.22962838366004
$ws.Range("P34").Value = 101.5955875616828
$ws.Range("R35").Value = 173.7492132756177
$ws.Range("L36").Value = 61.18167021676314
$ws.Range("N36").Value = 38.66169381481656
$ws.Range("K37").Value = 94.30397654773019
$ws.Range("J39").Value = 93.17061249236157
$ws.Range("N39").Value = 38.66169381481656
$ws.Range("P39").Value = 65.92768427608706
$ws.Range("Q39").Value = 94.49434172313325
$ws.Range("L40").Value = 90.4687457914608
$ws.Range("M40").Value = 92.09541281912071
$ws.Range("N40").Value = 81.96869489115805
$ws.Range("O40").Value = 96.22962838366004
$ws.Range("P40").Value = 101.5955875616828
$ws.Range("J41").Value = 124.5190384721106
$ws.Range("P41").Value = 135.4597561231036
$ws.Range("K42").Value = 80.29914934735042
$ws.Range("K43").Value = 94.30397654773019
$ws.Range("L43").Value = 90.4687457914608
$ws.Range("M43").Value = 92.09541281912071
$ws.Range("N43").Value = 81.96869489115805
$ws.Range("P43").Value = 0
$ws.Range("R44").Value = 173.7492132756177
$ws.Range("I45").Value = 87.25340171355576
$ws.Range("K46").Value = 0

# --- Household Surplus sheet ---
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B5").Value = 314358.8437207664
$ws.Range("B6").Value = 172190.7281481123
$ws.Range("B7").Value = 335250.0367186866
$ws.Range("B8").Value = 173168.8353292762
$ws.Range("B9").Value = 214113.0999756922
$ws.Range("B11").Value = 254347.4536560592
$ws.Range("B12").Value = 218907.6574858609
$ws.Range("B13").Value = 241219.9882795227
$ws.Range("B14").Value = 143835.8129211356
$ws.Range("B15").Value = 163843.3895748178
$ws.Range("B16").Value = 149997.7419296437

# --- Costs and Revenues sheet ---
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("E2").Value = 92547.98769918957
$ws.Range("F2").Value = 55990.47226622138
$ws.Range("G2").Value = 97920.00875579758
$ws.Range("H2").Value = 56241.98554137782
$ws.Range("I2").Value = 66770.51073617047
$ws.Range("K2").Value = 77116.48739683628
$ws.Range("L2").Value = 68003.39695307102
$ws.Range("M2").Value = 73740.85344286976
$ws.Range("N2").Value = 48699.20835071306
$ws.Range("O2").Value = 53844.01377594562
$ws.Range("P2").Value = 50283.70438147232
$ws.Range("E3").Value = 133100.0000000001
$ws.Range("B4").Value = 96756.66497830175
$ws.Range("C4").Value = 96756.66497830176
$ws.Range("D4").Value = 96756.66497830175
$ws.Range("E4").Value = 46997.99036755212
$ws.Range("F4").Value = 10440.47493458393
$ws.Range("G4").Value = 52370.01142416014
$ws.Range("H4").Value = 10691.98820974037
$ws.Range("I4").Value = 21220.51340453305
$ws.Range("K4").Value = 31566.49006519886
$ws.Range("L4").Value = 22453.39962143359
$ws.Range("M4").Value = 28190.85611123234
$ws.Range("N4").Value = 3149.211019075615
$ws.Range("O4").Value = 8294.016444308174
$ws.Range("P4").Value = 4733.707049834873
$ws.Range("B6").Value = -33627.59999999998
$ws.Range("C6").Value = -33627.6
$ws.Range("D6").Value = -33627.59999999998
$ws.Range("E6").Value = -91161.61477337092
$ws.Range("F6").Value = 41938.38522662912
$ws.Range("I6").Value = 41938.3852266291
$ws.Range("J6").Value = 41938.38522662912
$ws.Range("K6").Value = 41938.38522662909
$ws.Range("L6").Value = 41938.3852266291
$ws.Range("N6").Value = 41938.38522662912
$ws.Range("O6").Value = 41938.38522662912
